# Update code tinh luong % format cac bang
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LUY_KE_THANG_CAN_THO")

# Refresh last_edited_time (column D) for rows whose value was the
# 2024-07-17T17:23:00.000Z snapshot -> 2024-07-18T15:58:00.000Z
$newTimestamp = "2024-07-18T15:58:00.000Z"
$ws.Range("D4").Value = $newTimestamp
$ws.Range("D5").Value = $newTimestamp
$ws.Range("D6").Value = $newTimestamp
$ws.Range("D8").Value = $newTimestamp
$ws.Range("D12").Value = $newTimestamp
$ws.Range("D13").Value = $newTimestamp

# Update row 13 ("Thang 7") numeric properties
$ws.Range("W13").Value = 58922000    # properties.Chi tieu.number
$ws.Range("AA13").Value = 120358000  # properties.Luy ke.formula.number
$ws.Range("AE13").Value = 179280000  # properties.Tong doanh thu.formula.number
$ws.Range("AH13").Value = 153580000  # properties.Da thanh toan.number
$ws.Range("AK13").Value = 27         # properties.So luong don.number
$ws.Range("AN13").Value = 25700000   # properties.Thu no.number
$ws.Range("AQ13").Value = 179380000  # properties.Don gia.number
